# Reference steel factory from IEAGHG 2013
#
# Renames the existing IEAGHG_* steel process IDs to IEAGHGsteel_* (to
# disambiguate them from the newly-added shared Aux/Energy unit IDs), and
# appends four new unit-process rows: an auxiliary lime kiln, an auxiliary
# air separation unit, a one-step power plant, and a heat-recovery unit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename existing steel-process IDs (column A only; rows 9-14) ---
$ws.Cells.Item(9,  1).Value = "IEAGHGsteel_coke_oven"
$ws.Cells.Item(10, 1).Value = "IEAGHGsteel_sinter_plant"
$ws.Cells.Item(11, 1).Value = "IEAGHGsteel_blast_furnace"
$ws.Cells.Item(12, 1).Value = "IEAGHGsteel_BOF"
$ws.Cells.Item(13, 1).Value = "IEAGHGsteel_ladle"
$ws.Cells.Item(14, 1).Value = "IEAGHGsteel_forming"

# --- 2. Append new unit-process rows 15-18 ---
# Columns: A=ID, B=display name, C=product, D=productType,
#          E=varFile, F=varSheet, G=calcFile, H=calcSheet

$newRows = @(
    @{ Row=15; Vals=@("aux_lime kiln",      "Lime Kiln",      "CaO",            "outflow", "data/shared/AuxUnits_Variables.xlsx",    "Lime Kiln",      "data/shared/AuxUnits_Relationships.xlsx",    "Lime Kiln");      TextCols=@(1,2,3,4,5,6,7,8) },
    @{ Row=16; Vals=@("aux_air separation", "Air Separation", "O2",             "outflow", "data/shared/AuxUnits_Variables.xlsx",    "Air Separation", "data/shared/AuxUnits_Relationships.xlsx",    "Air Separation"); TextCols=@(1,2,3,4,5,6,7,8) },
    @{ Row=17; Vals=@("electricity_1step",  "Power Plant",    "electricity",    "outflow", "data/shared/EnergyUnits_Variables.xlsx", "one-step power", "data/shared/EnergyUnits_Relationships.xlsx", "one-step power"); TextCols=@(1,2,3,4,5,7) },
    @{ Row=18; Vals=@("heat_collector",     "heat recovery",  "recovered heat", "inflow",  "data/shared/EnergyUnits_Variables.xlsx", "Heat Recovery",  "data/shared/EnergyUnits_Relationships.xlsx", "Heat Recovery");  TextCols=@(1,2,3,4,5,7) }
)

foreach ($rowDef in $newRows) {
    $r = $rowDef.Row
    $vals = $rowDef.Vals
    $textCols = $rowDef.TextCols
    for ($c = 1; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($textCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $vals[$c - 1]
    }
}

# --- 3. Update the view: scroll back to top-left and move the selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("D18").Select() | Out-Null
